$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated frequency-table values (row 2)
$ws.Range("B2").Value = 0.0077720207253886
$ws.Range("C2").Value = 0.958549222797927
$ws.Range("D2").Value = 0.0466321243523316
$ws.Range("E2").Value = 0.987046632124352
$ws.Range("F2").Value = 0.989637305699482
$ws.Range("G2").Value = 0.0077720207253886
$ws.Range("H2").Value = 0.987046632124352
$ws.Range("I2").Value = 0.0310880829015544
$ws.Range("J2").Value = 0.955958549222798
$ws.Range("K2").Value = 0.00259067357512953
$ws.Range("L2").Value = 0.00518134715025907
$ws.Range("M2").Value = 0.00259067357512953
$ws.Range("O2").Value = 0.00518134715025907
$ws.Range("P2").Value = 0.00259067357512953
$ws.Range("U2").Value = 0.0129533678756477
$ws.Range("V2").Value = 0.00518134715025907
$ws.Range("W2").Value = 0.00259067357512953
$ws.Range("X2").Value = 0.00518134715025907

# Updated frequency-table values (row 3)
$ws.Range("B3").Value = 0.0077720207253886
$ws.Range("D3").Value = 0.953367875647668
$ws.Range("E3").Value = 0.0077720207253886
$ws.Range("H3").Value = 0.00259067357512953
$ws.Range("I3").Value = 0.132124352331606
$ws.Range("K3").Value = 0.00259067357512953
$ws.Range("L3").Value = 0.893782383419689
$ws.Range("M3").Value = 0.994818652849741
$ws.Range("N3").Value = 0.00518134715025907
$ws.Range("O3").Value = 0.0077720207253886
$ws.Range("R3").Value = 0.10880829015544
$ws.Range("S3").Value = 0.813471502590674
$ws.Range("V3").Value = 0.00518134715025907
$ws.Range("W3").Value = 0.00259067357512953
$ws.Range("X3").Value = 0.0077720207253886

# Updated frequency-table values (row 4)
$ws.Range("B4").Value = 0.981865284974093
$ws.Range("C4").Value = 0.0310880829015544
$ws.Range("F4").Value = 0.0077720207253886
$ws.Range("G4").Value = 0.989637305699482
$ws.Range("H4").Value = 0.0077720207253886
$ws.Range("I4").Value = 0.181347150259067
$ws.Range("J4").Value = 0.038860103626943
$ws.Range("K4").Value = 0.0077720207253886
$ws.Range("L4").Value = 0.00259067357512953
$ws.Range("N4").Value = 0.00259067357512953
$ws.Range("P4").Value = 0.99740932642487
$ws.Range("R4").Value = 0.00259067357512953
$ws.Range("U4").Value = 0.979274611398964
$ws.Range("V4").Value = 0.00518134715025907
$ws.Range("W4").Value = 0.994818652849741
$ws.Range("X4").Value = 0.984455958549223

# Updated frequency-table values (row 5)
$ws.Range("B5").Value = 0.00259067357512953
$ws.Range("C5").Value = 0.0077720207253886
$ws.Range("E5").Value = 0.00518134715025907
$ws.Range("F5").Value = 0.00259067357512953
$ws.Range("G5").Value = 0.00259067357512953
$ws.Range("H5").Value = 0.00259067357512953
$ws.Range("I5").Value = 0.655440414507772
$ws.Range("J5").Value = 0.00518134715025907
$ws.Range("K5").Value = 0.987046632124352
$ws.Range("L5").Value = 0.0984455958549223
$ws.Range("M5").Value = 0.00259067357512953
$ws.Range("N5").Value = 0.992227979274611
$ws.Range("O5").Value = 0.987046632124352
$ws.Range("R5").Value = 0.88860103626943
$ws.Range("S5").Value = 0.186528497409326
$ws.Range("U5").Value = 0.0077720207253886
$ws.Range("V5").Value = 0.984455958549223
$ws.Range("X5").Value = 0.00259067357512953
